$d = $word.ActiveDocument
$t = $d.Tables(1)
$values = @(
    "76-29=47",
    "7-3=4",
    "15+32=47",
    "27+17=44",
    "99-6=93",
    "50-7=43",
    "55-32=23",
    "44-15=29",
    "81-11=70",
    "57-46=11",
    "47+40=87",
    "13+73=86",
    "80-76=4",
    "0+22=22",
    "26+42=68",
    "57-49=8",
    "69-18=51",
    "54-15=39",
    "50-31=19",
    "75-36=39",
    "93-40=53",
    "43+47=90",
    "46+45=91",
    "60-46=14",
    "23+65=88",
    "65-1=64",
    "88-56=32",
    "8+11=19",
    "2+7=9",
    "22+58=80",
    "3+24=27",
    "67-28=39",
    "11+68=79",
    "57-26=31",
    "61+13=74",
    "42+15=57",
    "76+2=78",
    "38-5=33",
    "8+36=44",
    "2+48=50",
    "5+52=57",
    "47+46=93",
    "97-22=75",
    "10-4=6",
    "32-26=6",
    "15-8=7",
    "86-61=25",
    "24+74=98",
    "89+10=99",
    "0+40=40",
    "30-27=3",
    "2+5=7",
    "91-60=31",
    "61-38=23",
    "96+0=96",
    "78-65=13",
    "89+0=89",
    "61-34=27",
    "45+5=50",
    "54-36=18",
    "23-3=20",
    "92-86=6",
    "85-35=50",
    "35+35=70",
    "76-8=68",
    "25-21=4",
    "37-37=0",
    "60+6=66",
    "99-36=63",
    "34+21=55",
    "16+30=46",
    "63-26=37",
    "53-41=12",
    "35-31=4",
    "29+55=84",
    "6+87=93",
    "40+47=87",
    "78-20=58",
    "12+61=73",
    "37+59=96",
    "13+26=39",
    "19-11=8",
    "47+13=60",
    "98-53=45",
    "44+21=65",
    "56-34=22",
    "29+58=87",
    "95-27=68",
    "82-76=6",
    "93-84=9",
    "94-92=2",
    "53-45=8",
    "37-30=7",
    "90+2=92",
    "70-30=40",
    "26+66=92",
    "33+52=85",
    "14+4=18",
    "91-49=42",
    "71+23=94"
)
$count = 0
for ($row = 1; $row -le $t.Rows.Count; $row++) {
    for ($col = 1; $col -le $t.Columns.Count; $col++) {
        $c = $t.Cell($row, $col)
        $r = $c.Range
        $r.SetRange($r.Start, $r.End - 1)
        $r.Text = $values[$count]
        $count = $count + 1
    }
}
Write-Host "Updated cells:" $count
